$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 65: "loading sound {0}" -> "preloaded sound {0}" (text column D)
$ws.Cells.Item(65, 4).Value = "preloaded sound {0}"

# Row 91: new entries for cwl_log_stock_merge
$ws.Cells.Item(91, 1).Value = "cwl_log_stock_merge"
$ws.Cells.Item(91, 3).Value = "merged stock: {0} into character id: {1}"
$ws.Cells.Item(91, 4).Value = "merged stock: {0} into character id: {1}"

# Row 92: new entries for cwl_log_stock_add
$ws.Cells.Item(92, 1).Value = "cwl_log_stock_add"
$ws.Cells.Item(92, 3).Value = "added new stock: {0} to character id: {1}"
$ws.Cells.Item(92, 4).Value = "added new stock: {0} to character id: {1}"

# Update selection/view to match final workbook state
$ws.Application.Goto($ws.Range("A70"), $false)
$ws.Range("D91:D92").Select()
